$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "theta_threshold_range" parameter row (row 5) is removed entirely.
# Deleting the row shifts the old row 6 ("pie_threshold_range") up into row 5
# and drops the now-unused "theta_threshold_range" shared string.
$ws.Rows(5).Delete()

# Updated Min/Max values for the remaining parameters.
$ws.Range("B2").Value = 3.8
$ws.Range("C2").Value = 12.8

$ws.Range("B3").Value = 4.5
$ws.Range("C3").Value = 11.2

$ws.Range("B4").Value = 0.7
$ws.Range("C4").Value = 1.3

$ws.Range("B5").Value = 0
$ws.Range("C5").Value = 20

# Page setup (paper size / orientation) as set on the sheet.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Restore the active selection shown in the saved workbook.
$ws.Range("C3").Select()
